# Apply the recorded data update to the "Artfynd" sheet.
# All touched columns (Z, AB, D, F, G, H) hold plain text in this workbook
# (no number formats defined), so string literals are used for them to avoid
# Excel reinterpreting values such as "17:01" as a time serial number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 92180

# Row 3
$ws.Range("A3").Value = 131046847
$ws.Range("B3").Value = 79244
$ws.Range("E3").Value = 6425
$ws.Range("F3").Value = "Garnlav"
$ws.Range("G3").Value = "Alectoria sarmentosa"
$ws.Range("H3").Value = "(Ach.) Ach."
$ws.Range("Q3").Value = 402380
$ws.Range("R3").Value = 6818405
$ws.Range("Z3").Value = "17:01"
$ws.Range("AB3").Value = "17:01"

# Row 4
$ws.Range("A4").Value = 131047025
$ws.Range("B4").Value = 89194
$ws.Range("E4").Value = 510
$ws.Range("F4").Value = "Doftskinn"
$ws.Range("G4").Value = "Cystostereum murrayi"
$ws.Range("H4").Value = "(Berk. & M.A.Curtis.) Pouzar"
$ws.Range("Q4").Value = 402314
$ws.Range("R4").Value = 6818423
$ws.Range("Z4").Value = "16:05"
$ws.Range("AB4").Value = "16:05"

# Row 6
$ws.Range("B6").Value = 91809

# Row 8
$ws.Range("A8").Value = 131046843
$ws.Range("B8").Value = 79244
$ws.Range("Q8").Value = 402432
$ws.Range("R8").Value = 6818480
$ws.Range("Z8").Value = "16:20"
$ws.Range("AB8").Value = "16:20"

# Row 9
$ws.Range("A9").Value = 131046711
$ws.Range("B9").Value = 83224
$ws.Range("E9").Value = 6440
$ws.Range("F9").Value = "Vitgrynig nållav"
$ws.Range("G9").Value = "Chaenotheca subroscida"
$ws.Range("H9").Value = "(Eitner) Zahlbr."
$ws.Range("Q9").Value = 402363
$ws.Range("R9").Value = 6818428
$ws.Range("Z9").Value = "16:09"
$ws.Range("AB9").Value = "16:09"

# Row 10
$ws.Range("A10").Value = 131046844
$ws.Range("B10").Value = 79244
$ws.Range("E10").Value = 6425
$ws.Range("F10").Value = "Garnlav"
$ws.Range("G10").Value = "Alectoria sarmentosa"
$ws.Range("H10").Value = "(Ach.) Ach."
$ws.Range("Q10").Value = 402484
$ws.Range("R10").Value = 6818538
$ws.Range("Z10").Value = "16:23"
$ws.Range("AB10").Value = "16:23"

# Row 11
$ws.Range("B11").Value = 92268

# Row 14
$ws.Range("A14").Value = 131046806
$ws.Range("B14").Value = 83207
$ws.Range("D14").Value = "LC"
$ws.Range("E14").Value = 6439
$ws.Range("F14").Value = "Gulnål"
$ws.Range("G14").Value = "Chaenotheca brachypoda"
$ws.Range("H14").Value = "(Ach.) Tibell"
$ws.Range("Q14").Value = 402340
$ws.Range("R14").Value = 6818363
$ws.Range("Z14").Value = "17:05"
$ws.Range("AB14").Value = "17:05"

# Row 15
$ws.Range("A15").Value = 131046811
$ws.Range("B15").Value = 91829
$ws.Range("Q15").Value = 402450
$ws.Range("R15").Value = 6818298
$ws.Range("Z15").Value = "16:54"
$ws.Range("AB15").Value = "16:54"

# Row 16
$ws.Range("A16").Value = 131046808
$ws.Range("B16").Value = 91829
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 5432
$ws.Range("F16").Value = "Granticka"
$ws.Range("G16").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H16").Value = ""
$ws.Range("Q16").Value = 402323
$ws.Range("R16").Value = 6818416
$ws.Range("Z16").Value = "16:06"
$ws.Range("AB16").Value = "16:06"

# Row 18
$ws.Range("B18").Value = 91829

# Row 20
$ws.Range("B20").Value = 79244

# Row 22
$ws.Range("B22").Value = 91772

# Row 23
$ws.Range("B23").Value = 79244

# Row 25
$ws.Range("B25").Value = 79244
